$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 'Change the DB Table "Dissemination_reviews"'

$ws.Range("A12").Select()
